$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 586
$ws.Range("F5").Value = 2572
$ws.Range("F9").Value = 249
$ws.Range("F10").Value = 5302
$ws.Range("F12").Value = 1462
$ws.Range("F13").Value = 1379
$ws.Range("F14").Value = 593
$ws.Range("F15").Value = 6966
$ws.Range("F16").Value = 389
$ws.Range("F20").Value = 4670
$ws.Range("F23").Value = 2329
$ws.Range("F24").Value = 1255
$ws.Range("F25").Value = 442
$ws.Range("F26").Value = 1154
$ws.Range("F28").Value = 91
$ws.Range("F30").Value = 159
$ws.Range("F32").Value = 1273
$ws.Range("F34").Value = 237
$ws.Range("F35").Value = 519
$ws.Range("F36").Value = 202
$ws.Range("F37").Value = 1372
$ws.Range("F38").Value = 593
$ws.Range("F40").Value = 520
$ws.Range("F41").Value = 167
$ws.Range("F42").Value = 1117
$ws.Range("F43").Value = 2404
$ws.Range("F45").Value = 67
$ws.Range("F47").Value = 232

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 383
$ws.Range("F13").Value = 270
$ws.Range("F15").Value = 40
$ws.Range("F16").Value = 181
$ws.Range("F20").Value = 133
$ws.Range("F28").Value = 283

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 535
$ws.Range("F8").Value = 1291
$ws.Range("F10").Value = 1738
$ws.Range("F11").Value = 2183
$ws.Range("F12").Value = 614
$ws.Range("F13").Value = 517

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 586
$ws.Range("F6").Value = 535
$ws.Range("F7").Value = 2572
$ws.Range("F9").Value = 1291
$ws.Range("F10").Value = 2183
$ws.Range("F11").Value = 5302
$ws.Range("F12").Value = 614
$ws.Range("F17").Value = 1462
$ws.Range("F18").Value = 1379
$ws.Range("F19").Value = 593
$ws.Range("F20").Value = 6966
$ws.Range("F21").Value = 389
$ws.Range("F22").Value = 517
$ws.Range("F24").Value = 4670
$ws.Range("F25").Value = 2329
$ws.Range("F26").Value = 1255
$ws.Range("F27").Value = 442
$ws.Range("F28").Value = 1154
$ws.Range("F31").Value = 270
$ws.Range("F33").Value = 159
$ws.Range("F34").Value = 181
$ws.Range("F37").Value = 237
$ws.Range("F38").Value = 519
$ws.Range("F40").Value = 1372
$ws.Range("F42").Value = 167
$ws.Range("F44").Value = 1117
$ws.Range("F45").Value = 2404
$ws.Range("F46").Value = 67
$ws.Range("F47").Value = 232
